$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{ B2 = 973.9537847600009; E2 = 28982.37596598056; I2 = 16175.28135478; L2 = 48524.529503538; M2 = 10590.587968015; N2 = 7158.131594847755; O2 = 6981.145263461231 }
    "2030" = @{ B2 = 5712.560177842886; E2 = 56106.05588781912; I2 = 44217.8984721661; L2 = 66966.57749858923; M2 = 21984.28023276101; N2 = 10598.18910437709; O2 = 12064.7721182177 }
    "2035" = @{ A2 = 2861.961401238371; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15138.42652842583; O2 = 14762.98081419999 }
    "2040" = @{ A2 = 2861.961401238371; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15242.98365687085; O2 = 14762.98081419999 }
    "2045" = @{ A2 = 6302.873118834019; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15779.3841405625; O2 = 17100.41772974749 }
    "2050" = @{ A2 = 6302.873118834019; B2 = 8026.889663087295; E2 = 67297.73995507321; I2 = 59256.42575923612; L2 = 66966.57749858923; M2 = 25464.6214365565; N2 = 15779.3841405625; O2 = 17100.41772974749 }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $cellRefStr = [string]$cellRef
        $ws.Range($cellRefStr).Value = $cellValues[$cellRef]
    }
}
